$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 184; this pushes the existing rows 184-228
# down to 186-230 (and Excel auto-extends the used range / dimension).
$ws.Rows.Item(184).EntireRow.Insert()
$ws.Rows.Item(184).EntireRow.Insert()

# Row 184: new "Primera" quality record for 2022-04-12 (serial 44663)
$ws.Cells.Item(184, 1).Value = 11
$ws.Cells.Item(184, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(184, 3).Value = "Bíobío"
$ws.Cells.Item(184, 4).Value = 44663
$ws.Cells.Item(184, 5).Value = 8
$ws.Cells.Item(184, 6).Value = 100112008
$ws.Cells.Item(184, 7).Value = "Coliflor"
$ws.Cells.Item(184, 8).Value = "Sin especificar"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 2000
$ws.Cells.Item(184, 11).Value = 900
$ws.Cells.Item(184, 12).Value = 1000
$ws.Cells.Item(184, 13).Value = 950
$ws.Cells.Item(184, 14).Value = "$/unidad"
$ws.Cells.Item(184, 15).Value = "Región del Maule"
$ws.Cells.Item(184, 16).Value = 950
$ws.Cells.Item(184, 17).Value = 1
$ws.Cells.Item(184, 18).Value = "Hortaliza"

# Row 185: new "Segunda" quality record for 2022-04-12 (serial 44663)
$ws.Cells.Item(185, 1).Value = 11
$ws.Cells.Item(185, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(185, 3).Value = "Bíobío"
$ws.Cells.Item(185, 4).Value = 44663
$ws.Cells.Item(185, 5).Value = 8
$ws.Cells.Item(185, 6).Value = 100112008
$ws.Cells.Item(185, 7).Value = "Coliflor"
$ws.Cells.Item(185, 8).Value = "Sin especificar"
$ws.Cells.Item(185, 9).Value = "Segunda"
$ws.Cells.Item(185, 10).Value = 1000
$ws.Cells.Item(185, 11).Value = 800
$ws.Cells.Item(185, 12).Value = 800
$ws.Cells.Item(185, 13).Value = 800
$ws.Cells.Item(185, 14).Value = "$/unidad"
$ws.Cells.Item(185, 15).Value = "Región del Maule"
$ws.Cells.Item(185, 16).Value = 800
$ws.Cells.Item(185, 17).Value = 1
$ws.Cells.Item(185, 18).Value = "Hortaliza"
